$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.363.45"
$ws.Range("E2").Value = "  -3.61%  "
$ws.Range("D3").Value = "2.482.18"
$ws.Range("E3").Value = "  -4.70%  "
$ws.Range("E4").Value = "  +0.13%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "545.93"
$c.ClearFormats()
$ws.Range("E5").Value = "  -4.94%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "143.97"
$c.ClearFormats()
$ws.Range("E6").Value = "  -7.55%  "
$ws.Range("E7").Value = "  +0.01%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.608"
$c.ClearFormats()
$ws.Range("E8").Value = "  -2.54%  "
$ws.Range("D9").Value = "2.483.40"
$ws.Range("E9").Value = "  -4.56%  "
$ws.Range("E10").Value = "  -10.17%  "
$ws.Range("E11").Value = "  -1.85%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "5.36"
$c.ClearFormats()
$ws.Range("E12").Value = "  -8.68%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.352"
$c.ClearFormats()
$ws.Range("E13").Value = "  -7.34%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "25.87"
$c.ClearFormats()
$ws.Range("E14").Value = "  -7.76%  "
$ws.Range("D15").Value = "2.929.35"
$ws.Range("E15").Value = "  -4.64%  "
$ws.Range("D16").Value = "61.267.33"
$ws.Range("E16").Value = "  -3.53%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.0000163"
$c.ClearFormats()
$ws.Range("E17").Value = "  -8.78%  "
$ws.Range("D18").Value = "2.490.45"
$ws.Range("E18").Value = "  -3.86%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "11.05"
$c.ClearFormats()
$ws.Range("E19").Value = "  -8.05%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "6.95"
$c.ClearFormats()
$ws.Range("E20").Value = "  -7.76%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "4.16"
$c.ClearFormats()
$ws.Range("E21").Value = "  -8.26%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "319.35"
$c.ClearFormats()
$ws.Range("E22").Value = "  -6.78%  "
$ws.Range("E23").Value = "  +0.04%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "63.30"
$c.ClearFormats()
$ws.Range("E24").Value = "  -6.20%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "1.72"
$c.ClearFormats()
$ws.Range("E25").Value = "  -4.20%  "
$ws.Range("D26").Value = "0.0₃0989"
$ws.Range("E26").Value = "  -8.59%  "
$ws.Range("D27").Value = "2.618.38"
$ws.Range("E27").Value = "  -3.70%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.ClearFormats()
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.48"
$c.ClearFormats()
$ws.Range("E29").Value = "  -6.27%  "
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "531.43"
$c.ClearFormats()
$ws.Range("E30").Value = "  -8.40%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "8.30"
$c.ClearFormats()
$ws.Range("E31").Value = "  -9.40%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "7.62"
$c.ClearFormats()
$ws.Range("E32").Value = "  -3.38%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.149"
$c.ClearFormats()
$ws.Range("E33").Value = "  -7.56%  "
$ws.Range("E34").Value = "  -8.25%  "
$ws.Range("E35").Value = "  -9.90%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "5.80"
$c.ClearFormats()
$ws.Range("E36").Value = "  -11.22%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "4.85"
$c.ClearFormats()
$ws.Range("E37").Value = "  -9.20%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.ClearFormats()
$ws.Range("E38").Value = "  -0.06%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.375"
$c.ClearFormats()
$ws.Range("E39").Value = "  -6.63%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "18.34"
$c.ClearFormats()
$ws.Range("E40").Value = "  -7.01%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "143.67"
$c.ClearFormats()
$ws.Range("E41").Value = "  -6.95%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("E43").Value = "  -9.43%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "40.19"
$c.ClearFormats()
$ws.Range("E44").Value = "  -2.61%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "2.29"
$c.ClearFormats()
$ws.Range("E45").Value = "  -6.73%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "147.85"
$c.ClearFormats()
$ws.Range("E46").Value = "  -6.05%  "
$ws.Range("E47").Value = "  -8.64%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "20.83"
$c.ClearFormats()
$ws.Range("E48").Value = "  -11.25%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0533"
$c.ClearFormats()
$ws.Range("E49").Value = "  -9.20%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.585"
$c.ClearFormats()
$ws.Range("E50").Value = "  -6.91%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0937"
$c.ClearFormats()
$ws.Range("E51").Value = "  -6.52%  "
